# Applies the cryptos-list price/volume refresh described by the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.296.32"
$ws.Range("E2").Value = '  +2.80%  '

$ws.Range("D3").Value = "'3.813.40"
$ws.Range("E3").Value = '  +1.47%  '

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").Value = "'600.87"
$ws.Range("E5").Value = '  +0.90%  '

$ws.Range("D6").Value = "'171.04"

$ws.Range("D7").Value = "'3.811.63"
$ws.Range("E7").Value = '  +1.46%  '

$ws.Range("E8").Value = '  -0.08%  '

$ws.Range("E9").Value = '  -0.01%  '

$ws.Range("E10").Value = '  -1.69%  '

$ws.Range("E11").Value = '  +0.33%  '

$ws.Range("E12").Value = '  -0.34%  '

$ws.Range("D13").Value = "'0.0000265"
$ws.Range("E13").Value = '  -4.91%  '

$ws.Range("D14").Value = "'36.94"
$ws.Range("E14").Value = '  +0.66%  '

$ws.Range("D15").Value = "'4.452.83"
$ws.Range("E15").Value = '  +1.49%  '

$ws.Range("D16").Value = "'3.809.30"
$ws.Range("E16").Value = '  +1.28%  '

$ws.Range("D17").Value = "'69.225.58"
$ws.Range("E17").Value = '  +2.62%  '

$ws.Range("D18").Value = "'18.21"
$ws.Range("E18").Value = '  -2.61%  '

$ws.Range("D19").Value = "'7.09"
$ws.Range("E19").Value = '  -1.84%  '

$ws.Range("E20").Value = '  -0.24%  '

$ws.Range("D21").Value = "'11.04"
$ws.Range("E21").Value = '  +5.00%  '

$ws.Range("D22").Value = "'471.77"
$ws.Range("E22").Value = '  +0.52%  '

$ws.Range("D23").Value = "'0.710"
$ws.Range("E23").Value = '  -1.50%  '

$ws.Range("D24").Value = "'84.92"
$ws.Range("E24").Value = '  +1.20%  '

$ws.Range("E25").Value = '  +1.18%  '

$ws.Range("E26").Value = '  +0.65%  '

$ws.Range("D27").Value = "'12.25"
$ws.Range("E27").Value = '  +0.61%  '

$ws.Range("D28").Value = "'10.30"
$ws.Range("E28").Value = '  -1.48%  '

$ws.Range("E29").Value = '  +0.11%  '

$ws.Range("D30").Value = "'3.961.77"
$ws.Range("E30").Value = '  +1.43%  '

$ws.Range("E31").Value = '  -2.45%  '

$ws.Range("D32").Value = "'7.50"
$ws.Range("E32").Value = '  -2.43%  '

$ws.Range("D33").Value = "'2.25"
$ws.Range("E33").Value = '  +0.37%  '

$ws.Range("D34").Value = "'30.39"
$ws.Range("E34").Value = '  -0.74%  '

$ws.Range("D35").Value = "'9.44"
$ws.Range("E35").Value = '  +3.25%  '

$ws.Range("D37").Value = "'3.767.71"
$ws.Range("E37").Value = '  +1.25%  '

$ws.Range("E39").Value = '  -7.29%  '

$ws.Range("E40").Value = '  +1.33%  '

$ws.Range("E41").Value = '  +1.29%  '

$ws.Range("D42").Value = "'5.89"
$ws.Range("E42").Value = '  +0.32%  '

$ws.Range("D43").Value = "'0.998"
$ws.Range("E43").Value = '  -0.16%  '

$ws.Range("D44").Value = "'0.311"
$ws.Range("E44").Value = '  -0.49%  '

$ws.Range("E45").Value = '  +0.02%  '

$ws.Range("D46").Value = "'1.98"
$ws.Range("E46").Value = '  +1.42%  '

$ws.Range("B47").Value = 'Arweave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D47").Value = "'43.75"
$ws.Range("E47").Value = '  +10.71%  '

$ws.Range("B48").Value = 'Cosmos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D48").Value = "'8.66"
$ws.Range("E48").Value = '  -0.86%  '

$ws.Range("D49").Value = "'46.13"
$ws.Range("E49").Value = '  +0.60%  '

$ws.Range("D50").Value = "'402.81"
$ws.Range("E50").Value = '  +0.59%  '

$ws.Range("D51").Value = "'145.15"
$ws.Range("E51").Value = '  +3.47%  '
